$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.250.63"
$ws.Range("E2").Value = "  +1.15%  "
$ws.Range("D3").Value = "3.877.53"
$ws.Range("E3").Value = "  +0.85%  "
$ws.Range("E4").Value = "  -0.01%  "
$__style = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "471.83"
$ws.Range("D5").Style = $__style
$ws.Range("E5").Value = "  +10.72%  "
$__style = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.17"
$ws.Range("D6").Style = $__style
$ws.Range("E6").Value = "  +10.35%  "
$ws.Range("E7").Value = "  +3.36%  "
$__style = $ws.Range("D8").Style
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("D8").Style = $__style
$ws.Range("E8").Value = "  -0.05%  "
$__style = $ws.Range("D9").Style
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.745"
$ws.Range("D9").Style = $__style
$ws.Range("E9").Value = "  +1.94%  "
$__style = $ws.Range("D10").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.156"
$ws.Range("D10").Style = $__style
$ws.Range("E10").Value = "  -0.44%  "
$__style = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0000315"
$ws.Range("D11").Style = $__style
$ws.Range("E11").Value = "  -6.98%  "
$__style = $ws.Range("D12").Style
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "43.51"
$ws.Range("D12").Style = $__style
$ws.Range("E12").Value = "  +4.28%  "
$__style = $ws.Range("D13").Style
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "10.42"
$ws.Range("D13").Style = $__style
$ws.Range("E13").Value = "  -0.74%  "
$ws.Range("D14").Value = "4.518.25"
$ws.Range("E14").Value = "  +1.23%  "
$__style = $ws.Range("D15").Style
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.81"
$ws.Range("D15").Style = $__style
$ws.Range("E15").Value = "  -5.91%  "
$ws.Range("D16").Value = "3.886.76"
$ws.Range("E16").Value = "  +0.81%  "
$__style = $ws.Range("D18").Style
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "20.09"
$ws.Range("D18").Style = $__style
$ws.Range("E18").Value = "  +0.06%  "
$__style = $ws.Range("D19").Style
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.16"
$ws.Range("D19").Style = $__style
$ws.Range("E19").Value = "  +6.12%  "
$ws.Range("D20").Value = "67.532.70"
$ws.Range("E20").Value = "  +1.15%  "
$__style = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "435.16"
$ws.Range("D21").Style = $__style
$ws.Range("E21").Value = "  +4.83%  "
$__style = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "14.84"
$ws.Range("D22").Style = $__style
$ws.Range("E22").Value = "  -1.63%  "
$__style = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "89.07"
$ws.Range("D23").Style = $__style
$ws.Range("E23").Value = "  +5.02%  "
$ws.Range("E24").Value = "  +6.13%  "
$ws.Range("E25").Value = "  +11.06%  "
$__style = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "37.95"
$ws.Range("D26").Style = $__style
$ws.Range("E26").Value = "  +1.00%  "
$__style = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.27"
$ws.Range("D27").Style = $__style
$ws.Range("E27").Value = "  +11.61%  "
$__style = $ws.Range("D28").Style
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.95"
$ws.Range("D28").Style = $__style
$ws.Range("E28").Value = "  -0.62%  "
$__style = $ws.Range("D29").Style
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.48"
$ws.Range("D29").Style = $__style
$ws.Range("E29").Value = "  +2.49%  "
$__style = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "729.57"
$ws.Range("D30").Style = $__style
$ws.Range("E30").Value = "  +1.47%  "
$ws.Range("B31").Value = "Cosmos"
$ws.Range("C31").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$__style = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "13.82"
$ws.Range("D31").Style = $__style
$ws.Range("E31").Value = "  +0.47%  "
$ws.Range("B32").Value = "Hedera"
$ws.Range("C32").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$__style = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.135"
$ws.Range("D32").Style = $__style
$ws.Range("E32").Value = "  +7.99%  "
$__style = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.76"
$ws.Range("D33").Style = $__style
$ws.Range("E33").Value = "  -0.52%  "
$__style = $ws.Range("D34").Style
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "43.75"
$ws.Range("D34").Style = $__style
$ws.Range("E34").Value = "  +11.46%  "
$__style = $ws.Range("D35").Style
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.164"
$ws.Range("D35").Style = $__style
$ws.Range("E35").Value = "  +7.49%  "
$__style = $ws.Range("D36").Style
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "57.91"
$ws.Range("D36").Style = $__style
$ws.Range("E36").Value = "  +4.15%  "
$__style = $ws.Range("D37").Style
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.998"
$ws.Range("D37").Style = $__style
$ws.Range("E37").Value = "  -0.05%  "
$__style = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.52"
$ws.Range("D38").Style = $__style
$ws.Range("E38").Value = "  -5.80%  "
$__style = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0484"
$ws.Range("D39").Style = $__style
$ws.Range("E39").Value = "  +4.12%  "
$ws.Range("E40").Value = "  +10.21%  "
$ws.Range("B41").Value = "ThetaToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$__style = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.90"
$ws.Range("D41").Style = $__style
$ws.Range("E41").Value = "  -0.17%  "
$ws.Range("B42").Value = "PEPE"
$ws.Range("C42").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D42").Value = "0.0₃0688"
$ws.Range("E42").Value = "  -7.73%  "
$__style = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.142"
$ws.Range("D43").Style = $__style
$ws.Range("E43").Value = "  +3.99%  "
$__style = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.59"
$ws.Range("D44").Style = $__style
$ws.Range("E44").Value = "  +13.16%  "
$ws.Range("E45").Value = "  -0.09%  "
$__style = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.46"
$ws.Range("D46").Style = $__style
$ws.Range("E46").Value = "  +2.03%  "
$ws.Range("E47").Value = "  +7.11%  "
$ws.Range("E48").Value = "  +0.25%  "
$__style = $ws.Range("D49").Style
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.16"
$ws.Range("D49").Style = $__style
$ws.Range("E49").Value = "  +5.40%  "
$__style = $ws.Range("D50").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "144.45"
$ws.Range("D50").Style = $__style
$ws.Range("E50").Value = "  +1.83%  "
$ws.Range("E51").Value = "  +1.23%  "
